$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 325.7143
$ws.Range("I2").Value = 363.33334
$ws.Range("K2").Value = 363.33334
$ws.Range("M2").Value = -250.33334
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H17").Value = 578.5238000000001
$ws.Range("J17").Value = 671.6
$ws.Range("L17").Value = 2014.8
$ws.Range("N17").Value = -2350.8
$ws.Range("H33").Value = 262.7143
$ws.Range("I33").Value = 231.88889
$ws.Range("J33").Value = 318.2
$ws.Range("K33").Value = 231.88889
$ws.Range("L33").Value = 318.2
$ws.Range("M33").Value = -2.888890000000004
$ws.Range("N33").Value = -776.2
$ws.Range("H58").Value = 3316.5
$ws.Range("I58").Value = 107.5
$ws.Range("J58").Value = 4118.75
$ws.Range("K58").Value = 322.5
$ws.Range("L58").Value = 12356.25
$ws.Range("M58").Value = -172.5
$ws.Range("N58").Value = -12656.25
$ws.Range("H76").Value = 1900
$ws.Range("I76").Value = 1900
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 1900
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -1585
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 1900
$ws.Range("I79").Value = 1900
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 1900
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -808
$ws.Range("N79").ClearContents()
$ws.Range("H103").Value = 4159.8
$ws.Range("I103").Value = 2899.5
$ws.Range("K103").Value = 8698.5
$ws.Range("M103").Value = -8112.5

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 15000
$ws.Range("J24").Value = 15000
$ws.Range("L24").Value = 15000
$ws.Range("N24").Value = -15748
$ws.Range("H61").Value = 950.36365
$ws.Range("I61").Value = 845.5
$ws.Range("K61").Value = 845.5
$ws.Range("M61").Value = -633.5
$ws.Range("H97").Value = 536.5
$ws.Range("I97").Value = 513.7646999999999
$ws.Range("J97").Value = 665.3333
$ws.Range("K97").Value = 513.7646999999999
$ws.Range("L97").Value = 665.3333
$ws.Range("M97").Value = -17.76469999999995
$ws.Range("N97").Value = -1657.3333
$ws.Range("H100").Value = 15000
$ws.Range("J100").Value = 15000
$ws.Range("L100").Value = 15000
$ws.Range("N100").Value = -17164
$ws.Range("H136").Value = 950.36365
$ws.Range("I136").Value = 845.5
$ws.Range("K136").Value = 2536.5
$ws.Range("M136").Value = 13.5

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4000
$ws.Range("I20").Value = 4000
$ws.Range("K20").Value = 4000
$ws.Range("M20").Value = -3753
$ws.Range("H86").Value = 2177.647
$ws.Range("I86").Value = 2268.0667
$ws.Range("K86").Value = 2268.0667
$ws.Range("M86").Value = -1145.0667
$ws.Range("H89").Value = 2177.647
$ws.Range("I89").Value = 2268.0667
$ws.Range("K89").Value = 11340.3335
$ws.Range("M89").Value = -5724.333499999999

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1050
$ws.Range("I10").Value = 100
$ws.Range("J10").Value = 2000
$ws.Range("K10").Value = 100
$ws.Range("L10").Value = 2000
$ws.Range("M10").Value = 39
$ws.Range("N10").Value = -2278
$ws.Range("H50").Value = 16008.875
$ws.Range("I50").Value = 5041.5
$ws.Range("J50").Value = 19664.666
$ws.Range("K50").Value = 5041.5
$ws.Range("L50").Value = 19664.666
$ws.Range("M50").Value = -4416.5
$ws.Range("N50").Value = -20914.666
$ws.Range("H51").Value = 29997
$ws.Range("J51").Value = 29997
$ws.Range("L51").Value = 29997
$ws.Range("N51").Value = -31469
$ws.Range("H58").Value = 2308.5908
$ws.Range("I58").Value = 1204.5
$ws.Range("J58").Value = 5252.8335
$ws.Range("K58").Value = 1204.5
$ws.Range("L58").Value = 5252.8335
$ws.Range("M58").Value = -1001.5
$ws.Range("N58").Value = -5658.8335
$ws.Range("H59").Value = 53333.332
$ws.Range("J59").Value = 53333.332
$ws.Range("L59").Value = 53333.332
$ws.Range("N59").Value = -55623.332
$ws.Range("H60").Value = 10825.667
$ws.Range("J60").Value = 11996
$ws.Range("L60").Value = 11996
$ws.Range("N60").Value = -13018
$ws.Range("H61").Value = 29997
$ws.Range("J61").Value = 29997
$ws.Range("L61").Value = 29997
$ws.Range("N61").Value = -30693
$ws.Range("H69").Value = 24996.334
$ws.Range("I69").Value = 24996.334
$ws.Range("K69").Value = 24996.334
$ws.Range("M69").Value = -24247.334
$ws.Range("H72").Value = 24996.334
$ws.Range("I72").Value = 24996.334
$ws.Range("K72").Value = 74989.00199999999
$ws.Range("M72").Value = -71245.00199999999
$ws.Range("H99").Value = 13347.956
$ws.Range("I99").Value = 9085.583000000001
$ws.Range("K99").Value = 9085.583000000001
$ws.Range("M99").Value = -7587.583000000001
$ws.Range("H107").Value = 2712.7917
$ws.Range("I107").Value = 3340.6
$ws.Range("J107").Value = 1666.4445
$ws.Range("K107").Value = 3340.6
$ws.Range("L107").Value = 1666.4445
$ws.Range("M107").Value = -1420.6
$ws.Range("N107").Value = -5506.4445
$ws.Range("H126").Value = 13347.956
$ws.Range("I126").Value = 9085.583000000001
$ws.Range("K126").Value = 27256.749
$ws.Range("M126").Value = -24786.749
$ws.Range("H134").Value = 3322.9
$ws.Range("I134").Value = 2872.3333
$ws.Range("K134").Value = 8616.999899999999
$ws.Range("M134").Value = -6081.999899999999
$ws.Range("H136").Value = 2308.5908
$ws.Range("I136").Value = 1204.5
$ws.Range("J136").Value = 5252.8335
$ws.Range("K136").Value = 3613.5
$ws.Range("L136").Value = 15758.5005
$ws.Range("M136").Value = -1063.5
$ws.Range("N136").Value = -20858.5005

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 77199.16
$ws.Range("I2").Value = 200058
$ws.Range("K2").Value = 1200348
$ws.Range("M2").Value = -1200235
$ws.Range("H13").Value = 218.2
$ws.Range("I13").Value = 247.75
$ws.Range("K13").Value = 743.25
$ws.Range("M13").Value = -575.25
$ws.Range("H38").Value = 59.07143
$ws.Range("J38").Value = 76.59999999999999
$ws.Range("L38").Value = 229.8
$ws.Range("N38").Value = -923.8

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H70").Value = 9500
$ws.Range("I70").Value = 9500
$ws.Range("K70").Value = 9500
$ws.Range("M70").Value = -9230
$ws.Range("H73").Value = 9500
$ws.Range("I73").Value = 9500
$ws.Range("K73").Value = 9500
$ws.Range("M73").Value = -8564
$ws.Range("H132").Value = 1745.6
$ws.Range("I132").Value = 1448.1333
$ws.Range("K132").Value = 4344.3999
$ws.Range("M132").Value = -1814.3999

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1912
$ws.Range("I61").Value = 1349.5
$ws.Range("K61").Value = 1349.5
$ws.Range("M61").Value = -1147.5
$ws.Range("H113").Value = 1912
$ws.Range("I113").Value = 1349.5
$ws.Range("K113").Value = 1349.5
$ws.Range("M113").Value = 820.5
$ws.Range("H127").Value = 40000
$ws.Range("J127").Value = 40000
$ws.Range("L127").Value = 40000
$ws.Range("N127").Value = -49920

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()
$ws.Range("H20").Value = 42499
$ws.Range("J20").Value = 42499
$ws.Range("L20").Value = 42499
$ws.Range("N20").Value = -42979
$ws.Range("H132").Value = 1671.4375
$ws.Range("I132").Value = 1780.2307
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 5340.6921
$ws.Range("L132").Value = 3600
$ws.Range("M132").Value = -2810.6921
$ws.Range("N132").Value = -8660
$ws.Range("H133").Value = 60000
$ws.Range("J133").Value = 60000
$ws.Range("L133").Value = 60000
$ws.Range("N133").Value = -70120
$ws.Range("H136").Value = 2120.818
$ws.Range("I136").Value = 1047.4286
$ws.Range("K136").Value = 3142.2858
$ws.Range("M136").Value = -592.2857999999997
